$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 was a text/shared-string "ß"; change it to a plain number 42
$ws.Range("B8").Value = 42

# New row 10: sequential numbers 0..6 across B10:H10
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 6

# New row 11: formulas computing the same values via arithmetic
$ws.Range("B11").Formula = "=21-3*7"
$ws.Range("C11").Formula = "= 2 * 7 - 12-1"
$ws.Range("D11").Formula = "= 2 * 7 - 11 -1"
$ws.Range("E11").Formula = "= 2 * 7 - 12-1"
$ws.Range("F11").Formula = "= 2 * 7 - 13 -1"
$ws.Range("G11").Formula = "= 2 * 7 - 12-1"
$ws.Range("H11").Formula = "= 2 * 7 - 13 -1"

# Update selection to reflect the last active cell B11
$ws.Range("B11").Select()
